$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Separate the two script-file groups with a thin bottom border under row 3.
$ws.Range("A3:E3").Borders.Item(9).LineStyle = 1
$ws.Range("A3:E3").Borders.Item(9).Weight = 2

# New entry: SCRIPT/P01P04A/um0727.ssb (rows 568 and 571)
$ws.Range("C4").Value = " Project P is a dream I\'ve had\nfor many years…"
$ws.Range("C5").Value = " Someday, I know I\'ll make a huge\ndiscovery that will amaze everyone, all over\nthe world."
$ws.Range("A4").Value = "SCRIPT/P01P04A/um0727.ssb"
$ws.Range("D5").Value = " Когда-нибудь, я сделаю великое\nоткрытие, которое восхитит каждого в этом\nмире."
$ws.Range("D4").Value = " Проект П это то, о чём я\nмечтал долгие годы..."
$ws.Range("E4").Value = " Ðñïåëó Ð üóï óï, ï œæí ÿ\níåœóàì äïìãéå ãïäú..."
$ws.Range("E5").Value = " Ëïãäà-îéáôäû, ÿ òäåìàý âåìéëïå\nïóëñúóéå, ëïóïñïå âïòöéóéó ëàçäïãï â üóïí\níéñå."
$ws.Range("B4").Value = 568
$ws.Range("B5").Value = 571

$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(5).RowHeight = 31.8

# Matches the authored file's final cursor position.
$ws.Range("D4").Select() | Out-Null
